# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that a trailing "System" entry is moved to the front of the comma-separated
# list, preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7 ("Recorded By")
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -notlike "*,*") {
        continue
    }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -ceq "System") {
        $rest = $parts[0..($parts.Count - 2)]
        $newVal = "System, " + [string]::Join(", ", $rest)
        $cell.Value = $newVal
    }
}
